# Update stats with corrected data_correlations
$wb = $excel.ActiveWorkbook

# ---- Header renames (applied to both "Correlations" and "P Values" sheets) ----
$headers = @{
    "A1" = "Renal GSH (nmol/mg)";
    "B1" = "Renal GSSG (nmol/mg)";
    "C1" = "Renal Total Glutathione (nmol/mg)";
    "D1" = "Renal GSH/GSSG";
    "E1" = "Renal Eh (mV)";
    "G1" = "Hepatic GSH (nmol/mg)";
    "H1" = "Hepatic GSSG (nmol/mg)";
    "I1" = "Hepatic Total Glutathione (nmol/mg)";
    "J1" = "Hepatic GSH/GSSG";
    "K1" = "Hepatic Eh (mV)";
}

foreach ($sheetName in @("Correlations", "P Values")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $headers.Keys) {
        $ws.Range($addr).Value = $headers[$addr]
    }
}

# ---- Corrected correlation values on the "Correlations" sheet ----
$ws = $wb.Worksheets.Item("Correlations")

$data = @{
    2  = @("1", "0.3521064146258778", "0.9980528573398034", "0.4013309999581033", "-0.7732465092520669", "0.07121445498857891", "0.1989854281876762", "0.1568556770569907", "0.1993024726669089", "0.02221124464434825", "-0.1536022281255304");
    3  = @("0.3521064146258778", "1", "0.394263316972649", "-0.6606076551040662", "0.2625452042503621", "0.003355491492366705", "-0.05876753207235465", "0.01434459480465678", "-0.05602228799564749", "-0.1257951259286456", "0.1046137714338123");
    4  = @("0.9980528573398034", "0.394263316972649", "1", "0.3602168507545853", "-0.7432812758557322", "0.06704117126479425", "0.1927129152134672", "0.1544126655825603", "0.1931987257790364", "0.01469164365969326", "-0.1445993923683179");
    5  = @("0.4013309999581033", "-0.6606076551040662", "0.3602168507545853", "1", "-0.8790005084834919", "0.05935148474608827", "0.2185035389282103", "0.1124841865156832", "0.2158769791319677", "0.1497587271766938", "-0.2264377540702833");
    6  = @("-0.7732465092520669", "0.2625452042503621", "-0.7432812758557322", "-0.8790005084834919", "1", "-0.07916312637738215", "-0.2462622373808689", "-0.1472789699141999", "-0.2444318995834648", "-0.1241568421850889", "0.2376473570794778");
    7  = @("0.07121445498857891", "0.003355491492366705", "0.06704117126479425", "0.05935148474608827", "-0.07916312637738215", "1", "0.13369871861718", "0.09772014760320902", "0.1310647815837681", "0.03675699777175288", "-0.1230157056334216");
    8  = @("0.1989854281876762", "-0.05876753207235465", "0.1927129152134672", "0.2185035389282103", "-0.2462622373808689", "0.13369871861718", "1", "0.7905492366557961", "0.999682482891086", "0.1372934228397522", "-0.8100036720158353");
    9  = @("0.1568556770569907", "0.01434459480465678", "0.1544126655825603", "0.1124841865156832", "-0.1472789699141999", "0.09772014760320902", "0.7905492366557961", "1", "0.8014468546416662", "-0.4429971406708035", "-0.3256331147881447");
    10 = @("0.1993024726669089", "-0.05602228799564749", "0.1931987257790364", "0.2158769791319677", "-0.2444318995834648", "0.1310647815837681", "0.999682482891086", "0.8014468546416662", "1", "0.12005681992252", "-0.7989592333815935");
    11 = @("0.02221124464434825", "-0.1257951259286456", "0.01469164365969326", "0.1497587271766938", "-0.1241568421850889", "0.03675699777175288", "0.1372934228397522", "-0.4429971406708035", "0.12005681992252", "1", "-0.6524132569486403");
    12 = @("-0.1536022281255304", "0.1046137714338123", "-0.1445993923683179", "-0.2264377540702833", "0.2376473570794778", "-0.1230157056334216", "-0.8100036720158353", "-0.3256331147881447", "-0.7989592333815935", "-0.6524132569486403", "1");
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = [double]$values[$i]
    }
}

$wb.Save()
